# Add summary rows (averages / extremes) below the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the "k" column (J) used as the constraint value.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 12

# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").VerticalAlignment = -4108

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B15").Font.Bold = $true
$ws.Range("B15").VerticalAlignment = -4108

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B16").Font.Bold = $true
$ws.Range("B16").VerticalAlignment = -4108

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"
$ws.Range("B17").Font.Bold = $true
$ws.Range("B17").VerticalAlignment = -4108

# Selection / active cell as left by the author before saving.
$ws.Range("J12").Select()
